# Updated experiment result with new k for cross-validation
#
# Adds a "Row" header label to A1 (the blank corner cell above the
# model-name column) and refreshes the random_forest / lsboost result
# rows (B2:I3) with newly computed metric values. The old_model row
# (row 4) keeps its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label for column A
$ws.Range("A1").Value = "Row"

# random_forest row
$ws.Range("B2").Value = 3.2689362359390013
$ws.Range("C2").Value = 0.23053146938921021
$ws.Range("D2").Value = 2.6585772283272284
$ws.Range("E2").Value = 0.55202835654394378
$ws.Range("F2").Value = 0.74298610790777497
$ws.Range("G2").Value = 0.65546775846332062
$ws.Range("H2").Value = 0.44797164345605622
$ws.Range("I2").Value = 0.72959910978265552

# lsboost row
$ws.Range("B3").Value = 2.5189749571727784
$ws.Range("C3").Value = 0.17764280374984331
$ws.Range("D3").Value = 2.1973353129524527
$ws.Range("E3").Value = 0.32779036837539788
$ws.Range("F3").Value = 0.57252979693234995
$ws.Range("G3").Value = 0.54174933751293208
$ws.Range("H3").Value = 0.67220963162460212
$ws.Range("I3").Value = 0.82742848981409889
